$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData {
    param($rowA, $rowB)

    # Capture all F:V (columns 6-22) values from both rows before overwriting.
    $valsA = @()
    $valsB = @()
    for ($c = 6; $c -le 22; $c++) {
        $valsA += $ws.Cells.Item($rowA, $c).Value2
        $valsB += $ws.Cells.Item($rowB, $c).Value2
    }

    $i = 0
    for ($c = 6; $c -le 22; $c++) {
        $ws.Cells.Item($rowA, $c).Value2 = $valsB[$i]
        $ws.Cells.Item($rowB, $c).Value2 = $valsA[$i]
        $i = $i + 1
    }
}

# Swap the match data (columns F:V) between these row pairs, keeping
# A:E (Indice/pais/torneio/temporada/data_partida) fixed in place.
Swap-RowData 31 32
Swap-RowData 41 42
Swap-RowData 49 50

# Append three new match rows (64, 65, 66) at the end of the sheet.
# Copy formats from the last existing row (63) so the new rows inherit the
# same cell styles (bordered/bold Indice column, date-formatted data_partida).
$ws.Range("A63:V63").Copy()
$ws.Range("A64:V66").PasteSpecial(-4122)

$newRows = @(
    @{ Row = 64; A = 63; E = 45269.625; F = "Magra"; G = 1; H = "ASO Chlef"; I = 1;
       J = 2.38; K = "08/12/2023 15:43"; L = 2.42; M = "09/12/2023 14:55";
       N = 2.82; O = "08/12/2023 15:43"; P = 2.9;  Q = "09/12/2023 14:55";
       R = 3.47; S = "08/12/2023 15:43"; T = 3.35; U = "09/12/2023 14:55";
       V = "https://www.betexplorer.com/football/algeria/ligue-1/magra-aso-chlef/CjA5ex5g/" },
    @{ Row = 65; A = 64; E = 45269.625; F = "US Souf"; G = 0; H = "Saoura"; I = 1;
       J = 3.28; K = "08/12/2023 15:43"; L = 3.07; M = "09/12/2023 13:02";
       N = 2.91; O = "08/12/2023 15:43"; P = 2.92; Q = "09/12/2023 13:02";
       R = 2.41; S = "08/12/2023 15:43"; T = 2.57; U = "09/12/2023 14:31";
       V = "https://www.betexplorer.com/football/algeria/ligue-1/us-souf-saoura/YP0QjEJO/" },
    @{ Row = 66; A = 65; E = 45269.75;  F = "ES Setif"; G = 2; H = "Biskra"; I = 2;
       J = 1.69; K = "08/12/2023 15:43"; L = 1.48; M = "09/12/2023 17:56";
       N = 3.25; O = "08/12/2023 15:43"; P = 3.87; Q = "09/12/2023 17:57";
       R = 5.16; S = "08/12/2023 15:43"; T = 8.279999999999999; U = "09/12/2023 17:57";
       V = "https://www.betexplorer.com/football/algeria/ligue-1/es-setif-biskra/QgE9fdKa/" }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value2  = $r.A
    $ws.Cells.Item($row, 2).Value2  = "algeria"
    $ws.Cells.Item($row, 3).Value2  = "ligue-1"
    $ws.Cells.Item($row, 4).Value2  = "2023-2024"
    $ws.Cells.Item($row, 5).Value2  = $r.E
    $ws.Cells.Item($row, 6).Value2  = $r.F
    $ws.Cells.Item($row, 7).Value2  = $r.G
    $ws.Cells.Item($row, 8).Value2  = $r.H
    $ws.Cells.Item($row, 9).Value2  = $r.I
    $ws.Cells.Item($row, 10).Value2 = $r.J
    $ws.Cells.Item($row, 11).Value2 = $r.K
    $ws.Cells.Item($row, 12).Value2 = $r.L
    $ws.Cells.Item($row, 13).Value2 = $r.M
    $ws.Cells.Item($row, 14).Value2 = $r.N
    $ws.Cells.Item($row, 15).Value2 = $r.O
    $ws.Cells.Item($row, 16).Value2 = $r.P
    $ws.Cells.Item($row, 17).Value2 = $r.Q
    $ws.Cells.Item($row, 18).Value2 = $r.R
    $ws.Cells.Item($row, 19).Value2 = $r.S
    $ws.Cells.Item($row, 20).Value2 = $r.T
    $ws.Cells.Item($row, 21).Value2 = $r.U
    $ws.Cells.Item($row, 22).Value2 = $r.V
}
